# Leave Card update - 12/22/2023 10:59 AM
# Adds new leave credit/usage entries for Oct-Dec 2023, inserts a "2024" year-header
# row in the leave table, and appends two new monthly rows at the bottom so the
# table keeps its fixed monthly cadence.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1. Insert two blank rows right after row 211 (i.e. at 212:213), pushing the
#    remaining monthly rows (old 212..342) down to 214..344.
# ---------------------------------------------------------------------------
$ws.Rows("212:213").Insert() | Out-Null

# The plain row-insert does not fully preserve the table's cell formatting
# (borders get dropped), so re-apply formatting from row 211 which has the
# exact look of a normal/default table data row.
$ws.Range("A211:K211").Copy() | Out-Null
$ws.Range("A212:K213").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Restore the calculated "EARNED " helper-column formula for the two new rows.
$ws.Cells.Item(212, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Cells.Item(213, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Row 213 simply continues the monthly date sequence (Jan 1, 2024).
$ws.Cells.Item(213, 1).Value2 = 45292

# Row 212 becomes the "2024" year-header label (bold, quoted-text style), just
# like the existing "2023" header at row 196.
$ws.Cells.Item(212, 1).Value = "'2024"
$ws.Cells.Item(196, 1).Copy() | Out-Null
$ws.Cells.Item(212, 1).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Resize Table1 to cover the two new rows (A8:K344). Resizing auto-extends
#    the PERIOD date pattern into rows 343-344, but the auto-filled "EARNED "
#    helper formula uses a shorthand structured reference that this engine
#    mis-evaluates, so re-apply it in its fully-qualified form.
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A8:K344")) | Out-Null
$ws.Cells.Item(343, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Cells.Item(344, 7).Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# ---------------------------------------------------------------------------
# 3. Fill in the new leave entries on rows 208-211 (these rows keep their row
#    numbers; they previously held only a PERIOD date with blank details).
# ---------------------------------------------------------------------------

# Row 208: SL(1-0-0) earned 1.25, expiring 10/16/2023.
$ws.Cells.Item(207, 3).Copy() | Out-Null
$ws.Cells.Item(208, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(207, 11).Copy() | Out-Null
$ws.Cells.Item(208, 11).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(208, 2).Value2 = "SL(1-0-0)"
$ws.Cells.Item(208, 3).Value2 = 1.25
$ws.Cells.Item(208, 8).Value2 = 1
$ws.Cells.Item(208, 11).Value2 = 45215

# Row 209: SP(1-0-0) earned 1.25, expiring 12/1/2023.
$ws.Cells.Item(207, 3).Copy() | Out-Null
$ws.Cells.Item(209, 3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(207, 11).Copy() | Out-Null
$ws.Cells.Item(209, 11).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(209, 2).Value2 = "SP(1-0-0)"
$ws.Cells.Item(209, 3).Value2 = 1.25
$ws.Cells.Item(209, 11).Value2 = 45261

# Row 210: VL(1-0-0) used 1 day, dated 11/23/2023; PERIOD date cleared.
$ws.Cells.Item(207, 11).Copy() | Out-Null
$ws.Cells.Item(210, 11).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(210, 1).ClearContents() | Out-Null
$ws.Cells.Item(210, 2).Value2 = "VL(1-0-0)"
$ws.Cells.Item(210, 4).Value2 = 1
$ws.Cells.Item(210, 11).Value2 = 45253

# Row 211: VL(3-0-0) used 3 days, taken 12/27-29/2023; PERIOD date now 12/1/2023.
$ws.Cells.Item(211, 1).Value2 = 45261
$ws.Cells.Item(211, 2).Value2 = "VL(3-0-0)"
$ws.Cells.Item(211, 4).Value2 = 3
$ws.Cells.Item(211, 11).Value2 = "12/27-29/2023"

Write-Host "Leave card updated."
